$d = $word.ActiveDocument

# Locate the "7" in "Due: March 7th 2015" precisely via Find, so we don't
# depend on hard-coded character offsets.
$search = $d.Content
$found = $search.Find.Execute("March 7th", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'March 7th' text to correct the due date."
}

# "March " is 6 characters, so the "7" run starts right after it.
$sevenStart = $search.Start + 6
$sevenRange = $d.Range($sevenStart, $sevenStart + 1)
if ($sevenRange.Text -ne "7") {
    throw "Unexpected text at computed offset: [$($sevenRange.Text)]"
}

# Toggle Bold off/on around the text replacement. This forces the edited
# span to stay in its own run (matching the original run boundaries)
# instead of being silently coalesced into the neighbouring "March "/" "
# run that happens to share the same resolved formatting.
$sevenRange.Font.Bold = 0
$sevenRange.Text = "6"

$sixRange = $d.Range($sevenStart, $sevenStart + 1)
$sixRange.Font.Bold = 1

# Remove the old "_GoBack" bookmark (Word's last-edit-location marker) and
# re-add it collapsed right after the corrected "6", matching how Word
# drops _GoBack at the most recent edit point.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$gobackRange = $d.Range($sevenStart + 1, $sevenStart + 1)
$d.Bookmarks.Add("_GoBack", $gobackRange)
